# Add the "Script2" login credentials block to the Login sheet, below the
# existing (Script1) credentials block, following the same layout/style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# xlEdgeLeft=7, xlEdgeTop=8, xlEdgeBottom=9, xlEdgeRight=10
# xlInsideVertical=11, xlInsideHorizontal=12
# xlThin weight=2, xlMedium weight=-4138

function Format-HeaderRow($rng, $firstCell) {
    $rng.Font.Bold = $true
    $rng.Borders.Item(11).Weight = 2
    $rng.Borders.Item(7).Weight = -4138
    $rng.Borders.Item(10).Weight = -4138
    $rng.Borders.Item(8).Weight = -4138
    $rng.Borders.Item(9).Weight = -4138
    $firstCell.Font.Color = 255
}

function Format-DataRow($rng) {
    $rng.Borders.Item(7).Weight = 2
    $rng.Borders.Item(10).Weight = 2
    $rng.Borders.Item(9).Weight = 2
    $rng.Borders.Item(11).Weight = 2
}

# ---- Existing "Script1" block header/data already present (row 1-2); ----
# ---- give it the same box/bold formatting it has in the final file.   ----
$hdr1 = $ws.Range("A1:G1")
Format-HeaderRow $hdr1 $ws.Range("A1")

$data2 = $ws.Range("A2:G2")
Format-DataRow $data2

# Row 10: "Script2" section label - bold text on a yellow highlight.
$a10 = $ws.Range("A10")
$a10.Value = "Script2"
$a10.Font.Bold = $true
$a10.Interior.Color = 65535

# Row 11: header row for the Script2 block - identical labels/style to row 1.
$ws.Range("A11").Value = "Username"
$ws.Range("B11").Value = "password"
$ws.Range("C11").Value = "customerNumber"
$ws.Range("D11").Value = "coType"
$ws.Range("E11").Value = "PO"
$ws.Range("F11").Value = "ItemCode"
$ws.Range("G11").Value = "Quantity"

$hdr11 = $ws.Range("A11:G11")
Format-HeaderRow $hdr11 $ws.Range("A11")

# Row 12: data row for the Script2 block - new credentials.
$ws.Range("A12").Value = "AGSAutoT03"
$ws.Range("B12").Value = "SERVICE$08"
$ws.Range("C12").Value = "US00025065"
$ws.Range("D12").Value = "USA"
$ws.Range("E12:F12").NumberFormat = "@"
$ws.Range("E12").Value = "US-Gear-06"
$ws.Range("F12").Value = "TB7SX6CC"
$ws.Range("G12").Value = 1

$data12 = $ws.Range("A12:G12")
Format-DataRow $data12

# Row 3: blank spacer row, boxed with thin borders; F3 holds a stray
# leftover credential value (TB7SX1CC) typed with a text number format.
$row3 = $ws.Range("A3:G3")
$row3.Borders.LineStyle = 1
$row3.Borders.Weight = 2

$f3 = $ws.Range("F3")
$f3.NumberFormat = "@"
$f3.Value = "TB7SX1CC"

# Leave the selection where the user last clicked while editing.
$ws.Range("J11").Select()

Write-Host "Script2 credentials block added."
